$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a numeric-looking string to a cell while keeping it stored
# as text (matches the source data, which stores every Price/Volume cell as
# an inline string even when the text happens to look like a plain number).
# Flipping NumberFormat to "@" (Text) before the assignment stops Excel from
# parsing the literal into a Double; ClearFormats() afterwards drops the
# temporary formatting again so the cell keeps the workbook default style.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = '62.833.81'
$ws.Range("E2").Value = '  +2.07%  '

$ws.Range("D3").Value = '3.461.35'
$ws.Range("E3").Value = '  +2.00%  '

$ws.Range("E4").Value = '  -0.07%  '

Set-TextValue $ws.Range("D5") '577.54'
$ws.Range("E5").Value = '  +0.32%  '

Set-TextValue $ws.Range("D6") '147.50'
$ws.Range("E6").Value = '  +4.04%  '

$ws.Range("D7").Value = '3.463.23'
$ws.Range("E7").Value = '  +2.09%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  +1.69%  '

Set-TextValue $ws.Range("D10") '7.67'
$ws.Range("E10").Value = '  +0.38%  '

$ws.Range("E11").Value = '  +1.39%  '

$ws.Range("E12").Value = '  +4.13%  '

$ws.Range("D13").Value = '4.052.20'
$ws.Range("E13").Value = '  +1.96%  '

Set-TextValue $ws.Range("D14") '29.67'
$ws.Range("E14").Value = '  +6.12%  '

$ws.Range("E15").Value = '  +2.89%  '

$ws.Range("D16").Value = '3.464.92'
$ws.Range("E16").Value = '  +1.92%  '

$ws.Range("E17").Value = '  +0.45%  '

$ws.Range("D18").Value = '62.838.96'
$ws.Range("E18").Value = '  +1.94%  '

$ws.Range("E19").Value = '  +3.50%  '

Set-TextValue $ws.Range("D20") '14.31'
$ws.Range("E20").Value = '  +5.27%  '

Set-TextValue $ws.Range("D21") '9.18'
$ws.Range("E21").Value = '  +1.66%  '

Set-TextValue $ws.Range("D22") '388.64'
$ws.Range("E22").Value = '  +0.34%  '

Set-TextValue $ws.Range("D23") '0.556'
$ws.Range("E23").Value = '  +1.53%  '

Set-TextValue $ws.Range("D24") '74.44'
$ws.Range("E24").Value = '  -0.27%  '

$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("D26").Value = '3.603.74'
$ws.Range("E26").Value = '  +2.00%  '

$ws.Range("E27").Value = '  +1.49%  '

Set-TextValue $ws.Range("D28") '0.180'
$ws.Range("E28").Value = '  -7.04%  '

Set-TextValue $ws.Range("D29") '7.50'
$ws.Range("E29").Value = '  +1.36%  '

Set-TextValue $ws.Range("D30") '0.999'
$ws.Range("E30").Value = '  +0.06%  '

Set-TextValue $ws.Range("D31") '8.13'
$ws.Range("E31").Value = '  +1.94%  '

Set-TextValue $ws.Range("D32") '2.13'
$ws.Range("E32").Value = '  -0.84%  '

$ws.Range("E33").Value = '  +0.05%  '

$ws.Range("E34").Value = '  -2.33%  '

Set-TextValue $ws.Range("D35") '23.61'
$ws.Range("E35").Value = '  +1.38%  '

$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D36") '7.04'
$ws.Range("E36").Value = '  +1.91%  '

$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D37") '5.24'
$ws.Range("E37").Value = '  +3.88%  '

Set-TextValue $ws.Range("D38") '31.82'
$ws.Range("E38").Value = '  +19.84%  '

$ws.Range("E39").Value = '  +6.66%  '

Set-TextValue $ws.Range("D40") '169.91'
$ws.Range("E40").Value = '  +0.85%  '

$ws.Range("D41").Value = '3.499.31'
$ws.Range("E41").Value = '  +2.05%  '

Set-TextValue $ws.Range("D42") '0.0753'
$ws.Range("E42").Value = '  -1.33%  '

$ws.Range("E43").Value = '  +2.29%  '

Set-TextValue $ws.Range("D44") '42.37'
$ws.Range("E44").Value = '  -0.26%  '

Set-TextValue $ws.Range("D45") '4.46'
$ws.Range("E45").Value = '  +0.82%  '

$ws.Range("E46").Value = '  +3.21%  '

Set-TextValue $ws.Range("D47") '1.20'
$ws.Range("E47").Value = '  +3.89%  '

$ws.Range("D48").Value = '2.598.73'
$ws.Range("E48").Value = '  +5.80%  '

Set-TextValue $ws.Range("D49") '2.26'
$ws.Range("E49").Value = '  +11.57%  '

Set-TextValue $ws.Range("D50") '22.88'
$ws.Range("E50").Value = '  +1.00%  '

$ws.Range("E51").Value = '  +0.73%  '
